$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("stations")

# Insert a new row after the Samnanger_Tysseelva row (row 2), shifting the
# existing Sam_Stor / Sam_Froe / Sam_Frok rows down by one.
$ws.Rows.Item(3).Insert()

# The original "Samnanger_Tysseelva.shp" path has been superseded: the old
# shapefile is renamed/moved and a new "_ny" catchment shapefile is added
# for the new intercatchment station.
$ws.Cells.Item(2, 6).Value = "/home/jovyan/projects/critical_loads_2/cl_vestland/shapefiles/Samnanger.shp"

# New row: Samnanger_Tysseelva_Intercatchment station
$ws.Cells.Item(3, 1).Value = "Sam_Tyss_IntCat"
$ws.Cells.Item(3, 2).Value = "Samnanger_Tysseelva_Intercatchment"
$ws.Cells.Item(3, 4).Value = 5.7585797999999997
$ws.Cells.Item(3, 5).Value = 60.3748991
$ws.Cells.Item(3, 6).Value = "/home/jovyan/projects/critical_loads_2/cl_vestland/shapefiles/Samnanger_Tysseelva_ny.shp"

# Update the remaining shapefile paths to the new cl_vestland subfolder.
$ws.Cells.Item(4, 6).Value = "/home/jovyan/projects/critical_loads_2/cl_vestland/shapefiles/Samnanger_Storelva.shp"
$ws.Cells.Item(5, 6).Value = "/home/jovyan/projects/critical_loads_2/cl_vestland/shapefiles/Samnanger_Frolandselva.shp"
$ws.Cells.Item(6, 6).Value = "/home/jovyan/projects/critical_loads_2/cl_vestland/shapefiles/Samnanger_Frolandskanalen_ny.shp"

# Widen column F to fit the longer paths, and move the selection.
$ws.Columns.Item(6).ColumnWidth = 80.265625
$ws.Range("A4").Select()
